$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers: this runtime's COM layer always auto-converts a numeric-looking
# string assigned via .Value into a real number, and explicitly touching
# .NumberFormat materialises a brand-new style record (it never resolves
# back to an existing built-in style index). The only way found to land a
# literal (shared-string) text value on a cell while *reusing* an existing
# style index is: (1) force text entry with a leading apostrophe, then
# (2) copy ONLY the formatting from a cell that already has the desired
# style, via Copy + PasteSpecial(xlPasteFormats). The same Copy +
# PasteSpecial trick is used to move a numeric cell onto a different
# existing numeric style without creating a new style record.
# ---------------------------------------------------------------------------

function Set-TextStyled {
    param($ws, [string]$ref, [string]$text, [string]$styleFromRef)
    $target = $ws.Range($ref)
    $target.Value = "'" + $text
    $ws.Range($styleFromRef).Copy()
    $target.PasteSpecial(-4122) | Out-Null
}

function Set-NumStyled {
    param($ws, [string]$ref, $num, [string]$styleFromRef)
    $target = $ws.Range($ref)
    $target.Value = $num
    $ws.Range($styleFromRef).Copy()
    $target.PasteSpecial(-4122) | Out-Null
}

# Stable reference cells (untouched by this edit) carrying the two style
# indices ("s=14" general-text, "s=15" #,##0 number) we need to re-use.
$styleSrc14 = "A14"
$styleSrc15 = "I14"

# ---------------------------------------------------------------------------
# Header: volume/number and report week text (rich-text shared strings
# collapse to a single plain run; all runs shared identical formatting).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-TextStyled $ws "F14" "0" $styleSrc14
$ws.Range("N14").Value = 0

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-NumStyled $ws "C15" 1 $styleSrc15
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 250
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = 0

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 37
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = -21.276595744680
$ws.Range("L16").Value = 184.615384615385
$ws.Range("M16").Value = -19.565217391304
$ws.Range("N16").Value = -79.444444444444

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -85.714285714285
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -68.75
$ws.Range("I17").Value = 53
$ws.Range("J17").Value = 49
$ws.Range("K17").Value = 8.163265306122
$ws.Range("L17").Value = 43.243243243243
$ws.Range("M17").Value = 82.758620689655
$ws.Range("N17").Value = -20.895522388059

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 171
$ws.Range("J18").Value = 144
$ws.Range("K18").Value = 18.75
$ws.Range("L18").Value = 29.545454545454
$ws.Range("M18").Value = 35.714285714285
$ws.Range("N18").Value = -69.464285714285

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -35.294117647058
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -5.660377358490
$ws.Range("I19").Value = 359
$ws.Range("J19").Value = 314
$ws.Range("K19").Value = 14.331210191082
$ws.Range("L19").Value = 94.054054054054
$ws.Range("M19").Value = 90.957446808510
$ws.Range("N19").Value = 26.408450704225

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 6
Set-TextStyled $ws "D20" "0" $styleSrc14
Set-TextStyled $ws "E20" "***.*" $styleSrc14
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 81
$ws.Range("K20").Value = 65.306122448979
$ws.Range("L20").Value = 131.428571428571
$ws.Range("M20").Value = 6.578947368421
$ws.Range("N20").Value = -95.350172215843

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 94
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -7.843137254901
$ws.Range("I21").Value = 710
$ws.Range("J21").Value = 605
$ws.Range("K21").Value = 17.355371900826
$ws.Range("L21").Value = 76.178660049627
$ws.Range("M21").Value = 51.385927505330
$ws.Range("N21").Value = -75.017593244194

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -54.545454545454
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 43
$ws.Range("H24").Value = 25.581395348837
$ws.Range("I24").Value = 307
$ws.Range("J24").Value = 420
$ws.Range("K24").Value = -26.904761904761
$ws.Range("L24").Value = 9.252669039145
$ws.Range("M24").Value = 20.392156862745

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -52.631578947368
$ws.Range("I25").Value = 115
$ws.Range("J25").Value = 132
$ws.Range("K25").Value = -12.878787878787
$ws.Range("L25").Value = 64.285714285714
$ws.Range("M25").Value = 32.183908045977

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
Set-NumStyled $ws "C26" 1 $styleSrc15
$ws.Range("I26").Value = 8
$ws.Range("K26").Value = 166.666666666667
$ws.Range("L26").Value = 60

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-TextStyled $ws "D27" "0" $styleSrc14
Set-TextStyled $ws "E27" "***.*" $styleSrc14
$ws.Range("G27").Value = 1
$ws.Range("L27").Value = -27.272727272727

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-TextStyled $ws "D30" "0" $styleSrc14
Set-TextStyled $ws "E30" "***.*" $styleSrc14
